$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 137 - this shifts the existing rows 137:146 down to
# 138:147 and carries formatting down from the row above (matches the
# diff, which re-numbers the old rows 137-146 to 138-147 and adds a brand
# new row 137 with fresh weekly data).
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new weekly record. All the
# "constant" columns (A, B, C, E, F, G, H, I, N, O, Q, R) repeat the same
# values used throughout this block of rows.
$ws.Cells.Item(137, 1).Value = 3
$ws.Cells.Item(137, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(137, 3).Value = "Coquimbo"
$ws.Cells.Item(137, 4).Value = 44516
$ws.Cells.Item(137, 5).Value = 5
$ws.Cells.Item(137, 6).Value = 100112010
$ws.Cells.Item(137, 7).Value = "Achicoria"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 205
$ws.Cells.Item(137, 11).Value = 5000
$ws.Cells.Item(137, 12).Value = 6000
$ws.Cells.Item(137, 13).Value = 5476
$ws.Cells.Item(137, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(137, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(137, 16).Value = 342
$ws.Cells.Item(137, 17).Value = 16
$ws.Cells.Item(137, 18).Value = "Hortaliza"
